$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 87088
$ws.Range("D2").Value = 87230
$ws.Range("E2").Value = 86452
$ws.Range("F2").Value = 84312

$ws.Range("C3").Value = 18418
$ws.Range("D3").Value = 20323
$ws.Range("E3").Value = 16931
$ws.Range("F3").Value = 15986

$ws.Range("C4").Value = 14370
$ws.Range("D4").Value = 14755
$ws.Range("E4").Value = 13306
$ws.Range("F4").Value = 13075

$ws.Range("C5").Value = 18348
$ws.Range("D5").Value = 18030
$ws.Range("E5").Value = 17594
$ws.Range("F5").Value = 17515

$ws.Range("C6").Value = 18018
$ws.Range("D6").Value = 17456
$ws.Range("E6").Value = 15436
$ws.Range("F6").Value = 16389

$ws.Range("C7").Value = 263790
$ws.Range("D7").Value = 267335
$ws.Range("E7").Value = 67613
$ws.Range("F7").Value = 89633

$ws.Range("C8").Value = 60271
$ws.Range("D8").Value = 59706
$ws.Range("E8").Value = 58320
$ws.Range("F8").Value = 51432

$ws.Range("C9").Value = 14583
$ws.Range("D9").Value = 14777
$ws.Range("E9").Value = 15233
$ws.Range("F9").Value = 14196

$ws.Range("C10").Value = 75677
$ws.Range("D10").Value = 72284
$ws.Range("E10").Value = 63191
$ws.Range("F10").Value = 8458

$ws.Range("C11").Value = 39808
$ws.Range("D11").Value = 40924
$ws.Range("E11").Value = 38349
$ws.Range("F11").Value = 36973

$ws.Range("C12").Value = 21869
$ws.Range("D12").Value = 21523
$ws.Range("E12").Value = 20524
$ws.Range("F12").Value = 18963

$ws.Range("C13").Value = 196136
$ws.Range("D13").Value = 192244
$ws.Range("E13").Value = 46064
$ws.Range("F13").Value = 26318

$ws.Range("C14").Value = 13025
$ws.Range("D14").Value = 12466
$ws.Range("E14").Value = 12130
$ws.Range("F14").Value = 13584

$ws.Range("C15").Value = 88558
$ws.Range("D15").Value = 84975
$ws.Range("E15").Value = 81119
$ws.Range("F15").Value = 96687

$ws.Range("C16").Value = 61124
$ws.Range("D16").Value = 64212
$ws.Range("E16").Value = 53907
$ws.Range("F16").Value = 9616

$ws.Range("C17").Value = 200155
$ws.Range("D17").Value = 198974
$ws.Range("E17").Value = 184745
$ws.Range("F17").Value = 180124

$ws.Range("C18").Value = 30757
$ws.Range("D18").Value = 29415
$ws.Range("E18").Value = 25784
$ws.Range("F18").Value = 25129

$ws.Range("C19").Value = 11977
$ws.Range("D19").Value = 14003
$ws.Range("E19").Value = 8479
$ws.Range("F19").Value = 10332

$ws.Range("C20").Value = 7217
$ws.Range("D20").Value = 4000
$ws.Range("E20").Value = 4179
$ws.Range("F20").Value = 4235

$ws.Range("C21").Value = 10075
$ws.Range("D21").Value = 9265
$ws.Range("E21").Value = 9554
$ws.Range("F21").Value = 5293

$ws.Range("C22").Value = 47998
$ws.Range("D22").Value = 43715
$ws.Range("E22").Value = 42713
$ws.Range("F22").Value = 106247

$ws.Range("C23").Value = 107500
$ws.Range("D23").Value = 97644
$ws.Range("E23").Value = 108919
$ws.Range("F23").Value = 112338

$ws.Range("C24").Value = 68374
$ws.Range("D24").Value = 65476
$ws.Range("E24").Value = 52539
$ws.Range("F24").Value = 11967

$ws.Range("C25").Value = 92779
$ws.Range("D25").Value = 90320
$ws.Range("E25").Value = 84988
$ws.Range("F25").Value = 81503

$ws.Range("C26").Value = 108176
$ws.Range("D26").Value = 70778
$ws.Range("E26").Value = 80712
$ws.Range("F26").Value = 62093

$ws.Range("C27").Value = 28818
$ws.Range("D27").Value = 27365
$ws.Range("E27").Value = 26759
$ws.Range("F27").Value = 32838

$ws.Range("C28").Value = 37647
$ws.Range("D28").Value = 37201
$ws.Range("E28").Value = 35699
$ws.Range("F28").Value = 34183
